$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $result = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "NOT FOUND:" $old
    }
}

# 1. Updated object hash code in the MParagraphImpl reference
Replace-All "MParagraphImpl@19a20bb2" "MParagraphImpl@484b5a21"

# 2. Updated line numbers in the M2DocEvaluator / M2DocUtils / AbstractTemplatesTestSuite stack frames
Replace-All "M2DocEvaluator.caseQuery(M2DocEvaluator.java:561)" "M2DocEvaluator.caseQuery(M2DocEvaluator.java:586)"
Replace-All "M2DocEvaluator.java:1228)" "M2DocEvaluator.java:1239)"
Replace-All "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1437)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)"
Replace-All "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:288)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)"
Replace-All "M2DocEvaluator.generate(M2DocEvaluator.java:277)" "M2DocEvaluator.generate(M2DocEvaluator.java:281)"
Replace-All "M2DocUtils.generate(M2DocUtils.java:605)" "M2DocUtils.generate(M2DocUtils.java:805)"
Replace-All "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:461)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)"
Replace-All "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:368)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)"
Replace-All "GeneratedMethodAccessor76" "GeneratedMethodAccessor75"

# 3. Insert extra JUnit suite stack frames near the end of the trace.
#    Anchor on the unique tail sequence that ends the trace's inner block.
$tab = [char]9
$nl = [char]10

$oldTail = $tab + "at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl +
           $tab + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.run(ParentRunner.java:363)" + $nl +
           $tab + "at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)"

$newTail = $tab + "at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl +
           $tab + "at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)" + $nl +
           $tab + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.run(ParentRunner.java:363)" + $nl +
           $tab + "at org.junit.runners.Suite.runChild(Suite.java:128)" + $nl +
           $tab + "at org.junit.runners.Suite.runChild(Suite.java:27)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.run(ParentRunner.java:363)" + $nl +
           $tab + "at org.junit.runners.Suite.runChild(Suite.java:128)" + $nl +
           $tab + "at org.junit.runners.Suite.runChild(Suite.java:27)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)" + $nl +
           $tab + "at org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl +
           $tab + "at org.junit.runners.ParentRunner.run(ParentRunner.java:363)" + $nl +
           $tab + "at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)"

Replace-All $oldTail $newTail
